{"js": "// Move the \"Meta description\" paragraph from the top of the document down to\n// the end, turning it into two new-style paragraphs right before the final\n// (italic) image-prompt paragraph:\n//   1. A bold paragraph with the page title text.\n//   2. The review/meta text (without the \"Meta description\" label) replacing\n//      the old \"Create a cartoon-style...\" image-prompt paragraph, keeping\n//      the italic formatting already on that paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Meta description\" paragraph near the top of the doc.\nconst metaIndex = items.findIndex((p) => p.text.indexOf(\"Meta description\") !== -1);\nif (metaIndex === -1) {\n  throw new Error(\"Could not find the 'Meta description' paragraph\");\n}\nconst metaParagraph = items[metaIndex];\n\n// Locate the final paragraph (the italic AI image-prompt paragraph) - it is\n// simply the last paragraph in the body.\nconst lastParagraph = items[items.length - 1];\n\nconst titleText = \"Play Arabian Fire for Free - Unique and Exotic Slot Experience\";\nconst reviewText =\n  \"Read our review of Arabian Fire, the slot game with multiple jackpot opportunities, stunning visual design, and the chance to trigger the Loaded with Loot function. Play for free!\";\n\n// Insert a new bold paragraph right before the last paragraph with the title text.\nconst newTitleParagraph = lastParagraph.insertParagraph(titleText, \"Before\");\nnewTitleParagraph.font.bold = true;\n\n// Replace the last paragraph's text with the review text, keeping its\n// existing (italic) formatting.\nlastParagraph.insertText(reviewText, \"Replace\");\n\n// Remove the original \"Meta description\" paragraph from the top of the doc.\nmetaParagraph.delete();\n\nawait context.sync();\n", "ps1": "# Move the \"Meta description\" paragraph from the top of the document down to\n# the end, turning it into two paragraphs right before the final (italic)\n# image-prompt paragraph:\n#   1. A new bold paragraph with the page title text.\n#   2. The review/meta text (without the \"Meta description:\" label) replacing\n#      the text of the old \"Create a cartoon-style...\" image-prompt\n#      paragraph, while keeping that paragraph's existing italic formatting.\n\n$d = $word.ActiveDocument\n\n$titleText = \"Play Arabian Fire for Free - Unique and Exotic Slot Experience\"\n$reviewText = \"Read our review of Arabian Fire, the slot game with multiple jackpot opportunities, stunning visual design, and the chance to trigger the Loaded with Loot function. Play for free!\"\n\n# 1) Insert a new bold paragraph right before the last paragraph (the italic\n#    AI image-prompt paragraph), reusing the same run layout (\"<w:r/>\" empty\n#    run followed by a bold text run) that the rest of the document uses.\n#    Re-fetch the new (still-empty) paragraph by index rather than relying on\n#    a stale object reference, since $last re-seats onto the freshly\n#    inserted empty paragraph once InsertParagraphBefore() runs.\n$countBeforeInsert = $d.Paragraphs.Count\n$last = $d.Paragraphs.Last\n$last.Range.InsertParagraphBefore()\n$newPara = $d.Paragraphs.Item($countBeforeInsert)\n$xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $titleText + '</w:t></w:r></w:p>'\n$newPara.Range.InsertXML($xml)\n\n# 2) Replace the text of the final paragraph (the old \"Create a\n#    cartoon-style...\" image prompt) with the review text, excluding the\n#    trailing paragraph mark so the paragraph's own formatting is preserved.\n$lastRange = $d.Paragraphs.Last.Range\n$lastRange.MoveEnd(1, -1)\n$lastRange.Text = $reviewText\n\n# 3) Remove the original \"Meta description\" paragraph from the top of the\n#    document.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"Meta description*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
